$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.816.42"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.086.86"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.30"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.24"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "2.395.64"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.36"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "2.089.47"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "37.753.22"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.92"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.65"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0632"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0234"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.13"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0970"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "1.451.96"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "2.279.55"
$ws.Range("E51").Value = "  +0.19%  "
